# Insert a new "Skill Description" column after column A (SkillCode),
# containing a copy of the skill name from column A, and shift the
# existing SFIA Level / Keycode / Description columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row/column on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column before column B; this shifts B,C,D -> C,D,E.
$ws.Columns.Item(2).Insert()

# New header for column B.
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Populate column B for each data row with the skill name held in column A.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value()
}
